# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before column N, pushing the old N/O/P columns (Late /
# heading / Outstanding) one slot to the right (-> O/P/Q). The new
# column takes on the same width as the column to its left (M).
#
# The active sheet/tab also moves from "Summary" to "Repayment schedule",
# with a new selected cell on that sheet.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a fresh blank column before column N - this shifts the former
# N/O/P columns (and all their data/styles) one column to the right.
$wsRepay.Columns("N").Insert()

# The newly inserted column picks up the width that column M (the one
# now immediately to its left) already has.
$wsRepay.Columns("N").ColumnWidth = $wsRepay.Columns("M").ColumnWidth

# "Repayment schedule" becomes the active tab/sheet, with a new
# selection of S5 (previously the active sheet was "Summary", selection
# J3:J4).
$wsRepay.Activate()
$wsRepay.Range("S5").Select() | Out-Null
